# Auto-generated script to update cryptos worksheet values
# Mirrors the commit: "Updated cryptos list on Mon Feb 19 13:43:30 UTC 2024 with GitHub Actions"
#
# Price (column D) and Volume(1h) (column E) figures are refreshed for every
# coin row, two rows (Polygon/WrappedEther and Hedera/Cosmos) swapped rank
# position so their Coin/Link/Price/Volume columns are rewritten in place.
#
# NOTE: several Price values (e.g. "352.91") are numeric-looking strings that
# the workbook stores as TEXT (to match formats like "52.245.93" that are not
# valid numbers). Plain `.Value = "352.91"` assignment would be auto-coerced to
# a float by Excel, so for those cells we briefly force a Text number format,
# assign the value, then restore the cell style to "Normal" so no stray
# number-format/style is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range("D2").Value = "52.203.05"
$ws.Range("E2").Value = "  +1.14%  "

# Row 3
$ws.Range("D3").Value = "2.899.87"
$ws.Range("E3").Value = "  +3.83%  "

# Row 4
$ws.Range("E4").Value = "  +0.20%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "352.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.15%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "113.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.86%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.556"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.15%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.624"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.40%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.00%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0863"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.48%  "

# Row 12
$ws.Range("E12").Value = "  +0.59%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.28%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.05%  "

# Row 15
$ws.Range("D15").Value = "3.361.94"
$ws.Range("E15").Value = "  +4.02%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.915.75"
$ws.Range("E16").Value = "  +3.75%  "

# Row 17
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.995"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.81%  "

# Row 18
$ws.Range("D18").Value = "52.233.54"
$ws.Range("E18").Value = "  +1.27%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.85%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.45%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.86%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0976"
$ws.Range("E22").Value = "  +0.81%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.98%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.88%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.45%  "

# Row 26
$ws.Range("E26").Value = "  +8.50%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.67%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.05%  "

# Row 29
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.86%  "

# Row 30
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.103"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +15.81%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.64"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.92%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.47"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.71%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +12.07%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "53.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.07%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0449"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.98%  "

# Row 36
$ws.Range("E36").Value = "  -13.11%  "

# Row 37
$ws.Range("E37").Value = "  +0.00%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.64%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.35%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.05"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.31%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.52%  "

# Row 42
$ws.Range("E42").Value = "  +1.91%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.82%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.53%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "119.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.48%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.81%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.48%  "

# Row 48
$ws.Range("D48").Value = "2.176.64"
$ws.Range("E48").Value = "  +3.51%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.261"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +21.39%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0350"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +13.80%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.950"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.57%  "
